$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logout")

# Copy the formatting of the last existing row (row 8) down into the new
# row 9 so the new row picks up the same cell styling, then fill in the
# actual values for the new test case entry.
$ws.Range("A8:F8").Copy($ws.Range("A9:F9"))

$ws.Cells.Item(9, 1).Value = "testT2901"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = "Click"
$ws.Cells.Item(9, 5).Value = "Click"
$ws.Cells.Item(9, 6).Value = "Yes"

# Make the "Logout" sheet the active tab and select the new row's first
# cell, matching where the user was working.
$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
